# Update the default "owning group" participant values used by the
# Drools assignment rules on Sheet1 (column G, "Assign a literal
# participant" / ACTION column).
#
# The rules "Complaint - Default group", "Case File - Default group"
# and "DocumentRepository - Default group" assign the owning group to
# the ARKCASE_SUPERVISOR account, and the rules "Organization - Default
# group" and "Person - Default group" assign the owning group to the
# ARKCASE_ENTITY_ADMINISTRATOR account.
#
# Both accounts move from the old "000.<name>@APPDEV.ARMEDIA.COM" LDAP
# id to the new "<name>@ARMEDIA.COM" LDAP id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldSupervisor = "owning group, 000.ARKCASE_SUPERVISOR@APPDEV.ARMEDIA.COM"
$newSupervisor = "owning group, ARKCASE_SUPERVISOR@ARMEDIA.COM"

$oldEntityAdmin = "owning group, 000.ARKCASE_ENTITY_ADMINISTRATOR@APPDEV.ARMEDIA.COM"
$newEntityAdmin = "owning group, ARKCASE_ENTITY_ADMINISTRATOR@ARMEDIA.COM"

# The rule table lives in column G, rows 19-40. Walk every cell and
# replace whichever of the two old owning-group values is found, so the
# edit is resilient to the exact row positions. (Value2 is used to read
# the current contents since it returns plain text for string cells.)
for ($row = 19; $row -le 40; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2
    if ($current -eq $oldSupervisor) {
        $cell.Value = $newSupervisor
    }
    elseif ($current -eq $oldEntityAdmin) {
        $cell.Value = $newEntityAdmin
    }
}
